{"js": "// Update the intro paragraph and rebuild the car table (3 cols -> 4 cols)\n// with new header/separator/data rows, per the target diff.\n\nconst body = context.document.body;\n\n// 1) Update intro paragraph text (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst introPara = paragraphs.items[0];\nintroPara.load(\"text\");\nawait context.sync();\n\nintroPara.insertText(\n  \"Introducing a variety of car models across different brands, showcasing their engines and pricing points.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// 2) Locate the (only) table in the document.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// 3) Add the new 4th column (\"Price Range ($)\") at the end, seeded with\n//    the new price-range values for every existing row.\nconst newColumnValues = [\n  [\"Price Range ($) \"],\n  [\"nan\"],\n  [\"$80,000+        \"],\n  [\"$20,000-30,000  \"],\n  [\"$20,000-30,000  \"],\n  [\"$25,000-40,000  \"],\n  [\"$40,000+        \"],\n];\ntable.addColumns(\"End\", 1, newColumnValues);\nawait context.sync();\n\n// 4) Overwrite the first three columns with the new header/sep/data text.\nconst newFirstThreeColumns = [\n  [\"Car Model       \", \"Brand               \", \"Engine Type \"],\n  [\"------------------\", \"-------------\", \"---------------\"],\n  [\"Tesla Model S     \", \"Tesla         \", \"Electric        \"],\n  [\"Honda Civic       \", \"Honda          \", \"Gasoline        \"],\n  [\"Toyota Camry       \", \"Toyota         \", \"Gasoline        \"],\n  [\"Ford Mustang       \", \"Ford           \", \"Gasoline        \"],\n  [\"Chevrolet Bolt EV  \", \"Chevrolet      \", \"Electric        \"],\n];\n\nfor (let r = 0; r < newFirstThreeColumns.length; r++) {\n  for (let c = 0; c < newFirstThreeColumns[r].length; c++) {\n    table.getCell(r, c).value = newFirstThreeColumns[r][c];\n  }\n}\nawait context.sync();\n\n// 5) Resize every column to 2160 twips (= 108 points) \u2014 columnWidth is in\n//    points and rewrites the whole column (gridCol + every cell's tcW).\nfor (let c = 0; c < 4; c++) {\n  table.getCell(0, c).columnWidth = 108;\n}\nawait context.sync();\n", "ps1": "# Update the intro paragraph and rebuild the car table (3 cols -> 4 cols)\n# with new header/separator/data rows, per the target diff.\n\n$d = $word.ActiveDocument\n\n# 1) Update intro paragraph text (first paragraph in the body).\n$p = $d.Paragraphs.Item(1)\n$p.Range.Text = \"Introducing a variety of car models across different brands, showcasing their engines and pricing points.\"\n\n# 2) Locate the (only) table in the document.\n$t = $d.Tables.Item(1)\n\n# 3) Add a new 4th column at the end of the table.\n$t.Columns.Add() | Out-Null\n\n# 4) Fill in the new (4th) column's price-range values for every row.\n$priceCol = @(\n    \"Price Range ($) \",\n    \"nan\",\n    \"`$80,000+        \",\n    \"`$20,000-30,000  \",\n    \"`$20,000-30,000  \",\n    \"`$25,000-40,000  \",\n    \"`$40,000+        \"\n)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $t.Cell($r, 4).Range.Text = $priceCol[$r - 1]\n}\n\n# 5) Overwrite the first three columns with the new header/sep/data text.\n$newData = @(\n    @(\"Car Model       \", \"Brand               \", \"Engine Type \"),\n    @(\"------------------\", \"-------------\", \"---------------\"),\n    @(\"Tesla Model S     \", \"Tesla         \", \"Electric        \"),\n    @(\"Honda Civic       \", \"Honda          \", \"Gasoline        \"),\n    @(\"Toyota Camry       \", \"Toyota         \", \"Gasoline        \"),\n    @(\"Ford Mustang       \", \"Ford           \", \"Gasoline        \"),\n    @(\"Chevrolet Bolt EV  \", \"Chevrolet      \", \"Electric        \")\n)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowVals = $newData[$r - 1]\n    for ($c = 1; $c -le 3; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowVals[$c - 1]\n    }\n}\n\n# 6) Resize every column to 2160 twips (= 108 points).\nfor ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $t.Columns.Item($c).Width = 108\n}\n"}
